# "Generate Report for Archive"
# The localization status report is regenerated: every cell that used to
# read "Ready for handoff" now reads "In Translation" (the file moved from
# the handoff stage to being actively translated), on the Overview sheet as
# well as on each per-locale detail sheet. Because the new status text is
# shorter than the old one, the status columns that were sized to fit it
# are re-fitted to the new, narrower content.

$wb = $excel.ActiveWorkbook

# Update the status text everywhere it appears (Overview + per-locale sheets).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Re-fit the status columns now that the text is shorter. The new width
# matches what Excel computes when it re-fits those columns to the shorter
# "In Translation" text.
$newStatusColumnWidth = 12.5

#  - Overview sheet: the "zh-cn" (E) and "de-de" (F) status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns("E:F").ColumnWidth = $newStatusColumnWidth

#  - zh-cn detail sheet: the "Status" (C) column
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns("C:C").ColumnWidth = $newStatusColumnWidth

#  - de-de detail sheet: the "Status" (C) column
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns("C:C").ColumnWidth = $newStatusColumnWidth
